$d = $word.ActiveDocument
$t = $d.Tables.Item(1)

$updates = @(
    @{ Row = 1;  Col = 1; Text = "51×55=2805" },
    @{ Row = 1;  Col = 2; Text = "78×45=3510" },
    @{ Row = 1;  Col = 3; Text = "27×33=891" },
    @{ Row = 1;  Col = 4; Text = "88×90=7920" },
    @{ Row = 1;  Col = 5; Text = "95×97=9215" },

    @{ Row = 5;  Col = 1; Text = "28×45=1260" },
    @{ Row = 5;  Col = 2; Text = "30×65=1950" },
    @{ Row = 5;  Col = 3; Text = "26×88=2288" },
    @{ Row = 5;  Col = 4; Text = "90×55=4950" },
    @{ Row = 5;  Col = 5; Text = "38×22=836" },

    @{ Row = 10; Col = 1; Text = "46×69=3174" },
    @{ Row = 10; Col = 2; Text = "68×25=1700" },
    @{ Row = 10; Col = 3; Text = "22×61=1342" },
    @{ Row = 10; Col = 4; Text = "73×38=2774" },
    @{ Row = 10; Col = 5; Text = "65×24=1560" },

    @{ Row = 15; Col = 1; Text = "51×13=663" },
    @{ Row = 15; Col = 2; Text = "54×19=1026" },
    @{ Row = 15; Col = 3; Text = "77×89=6853" },
    @{ Row = 15; Col = 4; Text = "12×72=864" },
    @{ Row = 15; Col = 5; Text = "33×31=1023" },

    @{ Row = 20; Col = 1; Text = "78×90=7020" },
    @{ Row = 20; Col = 2; Text = "92×73=6716" },
    @{ Row = 20; Col = 3; Text = "93×67=6231" },
    @{ Row = 20; Col = 4; Text = "59×85=5015" },
    @{ Row = 20; Col = 5; Text = "75×38=2850" }
)

foreach ($u in $updates) {
    $cell = $t.Cell($u.Row, $u.Col)
    $cell.Range.Text = $u.Text
}

Write-Output "done"
